# Applies the diff: adds <w:noProof/> to the title's paragraph-mark rPr,
# adds <w:lang w:val="ru-RU"/> to the title's ": " run, and inserts a new
# paragraph (with a Judge-contest HYPERLINK field) right after the title.

$d = $word.ActiveDocument

# --- Step 1: replace paragraph 1 (title) with a version carrying the two added rPr children ---
$p1 = $d.Paragraphs(1)
$xmlPara1 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:p w14:paraId="1E462AC1" w14:textId="77777777" w:rsidR="00276520" w:rsidRPr="00276520" w:rsidRDefault="00276520" w:rsidP="00276520"><w:pPr><w:pStyle w:val="Heading1"/><w:jc w:val="center"/><w:rPr><w:noProof/><w:szCs w:val="40"/><w:lang w:val="bg-BG"/></w:rPr></w:pPr><w:r w:rsidRPr="00276520"><w:rPr><w:noProof/><w:szCs w:val="40"/><w:lang w:val="bg-BG"/></w:rPr><w:t>Упражнение</w:t></w:r><w:r w:rsidRPr="00276520"><w:rPr><w:szCs w:val="40"/><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r w:rsidRPr="00276520"><w:rPr><w:noProof/><w:szCs w:val="40"/><w:lang w:val="bg-BG"/></w:rPr><w:t>Моделиране на бази данни</w:t></w:r></w:p></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$p1.Range.InsertXML($xmlPara1)

# --- Step 2: split off a brand-new empty paragraph right after paragraph 1 ---
$p1 = $d.Paragraphs(1)
$p1.Range.InsertParagraphAfter()

# --- Step 3: fill that new paragraph with the Judge-link content (rStyle applied separately below) ---
$newPara = $d.Paragraphs(2)
$xmlPara2 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:p><w:pPr><w:rPr><w:lang w:val="bg-BG"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="bg-BG"/></w:rPr><w:t xml:space="preserve">Тествайте решенията си в </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>Judge</w:t></w:r><w:r><w:rPr><w:b/><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/><w:lang w:val="bg-BG"/></w:rPr><w:t>системата</w:t></w:r><w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r><w:rPr><w:lang w:val="bg-BG"/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:rPr><w:lang w:val="bg-BG"/></w:rPr><w:instrText>HYPERLINK "https://judge.softuni.org/Contests/Practice/Index/4606#0"</w:instrText></w:r><w:r><w:rPr><w:lang w:val="bg-BG"/></w:rPr></w:r><w:r><w:rPr><w:lang w:val="bg-BG"/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:lang w:val="bg-BG"/></w:rPr><w:t>https://judge.softuni.org/Contests/Practice/Index/4606#0</w:t></w:r><w:r><w:rPr><w:lang w:val="bg-BG"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r></w:p></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$newPara.Range.InsertXML($xmlPara2)

# --- Step 4: apply the Hyperlink character style to the visible URL text run ---
$newPara = $d.Paragraphs(2)
$linkRange = $newPara.Range.Duplicate
$found = $linkRange.Find.Execute("https://judge.softuni.org/Contests/Practice/Index/4606#0", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
  $linkRange.Style = "Hyperlink"
}

Write-Output "Edit complete. Paragraph count: $($d.Paragraphs.Count)"
